# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columns temporalidad (A), mes-nombre (H), modalidad (J),
# grupo-de-tipo-de-contrato (K) and sexo (M) move from being "dim"
# (dimension, skos:Concept, with an external mapping-*.xlsx file) to
# being "medida" (measure, xsd:int) and no longer reference a mapping
# workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "H", "J", "K", "M")

foreach ($col in $columns) {
    # Row 2: iaest-dimension:<x>  ->  iaest-measure:<x>
    $row2 = $ws.Range($col + "2")
    $row2.Value2 = ($row2.Value2 -replace "^iaest-dimension:", "iaest-measure:")

    # Row 3: dim -> medida
    $ws.Range($col + "3").Value2 = "medida"

    # Row 4: skos:Concept -> xsd:int
    $ws.Range($col + "4").Value2 = "xsd:int"

    # Row 5: the mapping-*.xlsx reference is removed entirely (cell itself
    # goes away, not just its value)
    $ws.Range($col + "5").Clear()
}
